$d = $word.ActiveDocument

# Step 1: Replace "Swetha git project" with itself using Find/Replace to
# normalize/simplify the runs and drop the spell-check proofErr markers.
$d.Content.Find.Execute("Swetha git project", $false, $false, $false, $false, $false, $true, 1, $false, "Swetha git project", 2) | Out-Null

# Step 2: Insert a new paragraph after the first paragraph, then an empty
# paragraph, per the target layout:
#   P1: Swetha git project
#   P2: (empty)
#   P3: Swetha brach is created  [+ bookmark _GoBack]
$para1 = $d.Paragraphs(1)
$rng = $para1.Range
$rng.Collapse(0)  # wdCollapseEnd
$rng.InsertParagraphAfter()

$para2 = $d.Paragraphs(2)
$rng2 = $para2.Range
$rng2.Collapse(0)
$rng2.InsertParagraphAfter()

$para3 = $d.Paragraphs(3)
$para3.Range.Text = "Swetha brach is created"

$d.Content.Find.Execute("Swetha brach is created", $false, $false, $false, $false, $false, $true, 1, $false, "Swetha brach is created", 2) | Out-Null
